$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.981.86"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.282.20"
$ws.Range("E3").Value = "  +1.57%  "
$rng = $ws.Range("D4")
$rng.Formula = "=""1.01"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = "  +0.14%  "
$rng = $ws.Range("D5")
$rng.Formula = "=""112.31"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  -2.53%  "
$rng = $ws.Range("D6")
$rng.Formula = "=""308.92"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  +6.48%  "
$rng = $ws.Range("D7")
$rng.Formula = "=""0.632"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  -0.22%  "
$rng = $ws.Range("D9")
$rng.Formula = "=""0.612"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  -0.91%  "
$rng = $ws.Range("D10")
$rng.Formula = "=""44.22"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  -4.95%  "
$rng = $ws.Range("D11")
$rng.Formula = "=""0.0925"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  -1.02%  "
$rng = $ws.Range("D12")
$rng.Formula = "=""55.11"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  -1.43%  "
$rng = $ws.Range("D13")
$rng.Formula = "=""8.79"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  -3.86%  "
$rng = $ws.Range("D14")
$rng.Formula = "=""1.09"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  +21.94%  "
$ws.Range("E15").Value = "  -0.36%  "
$rng = $ws.Range("D16")
$rng.Formula = "=""15.55"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "2.625.59"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "2.328.21"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").Value = "42.943.04"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  -0.51%  "
$rng = $ws.Range("D21")
$rng.Formula = "=""7.18"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  -3.79%  "
$rng = $ws.Range("D22")
$rng.Formula = "=""75.85"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  +2.69%  "
$rng = $ws.Range("D23")
$rng.Formula = "=""3.62"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Value = "  +4.91%  "
$ws.Range("E24").Value = "  +4.40%  "
$rng = $ws.Range("D25")
$rng.Formula = "=""255.27"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  +9.61%  "
$rng = $ws.Range("D26")
$rng.Formula = "=""8.95"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E26").Value = "  -3.87%  "
$rng = $ws.Range("D27")
$rng.Formula = "=""11.73"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  -3.74%  "
$rng = $ws.Range("D28")
$rng.Formula = "=""0.999"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = "  -0.10%  "
$rng = $ws.Range("D29")
$rng.Formula = "=""2.23"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -0.51%  "
$rng = $ws.Range("D30")
$rng.Formula = "=""38.16"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  -5.21%  "
$rng = $ws.Range("D31")
$rng.Formula = "=""174.45"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  -0.61%  "
$rng = $ws.Range("D32")
$rng.Formula = "=""22.12"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("E33").Value = "  -3.01%  "
$rng = $ws.Range("D34")
$rng.Formula = "=""0.0899"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  -1.74%  "
$rng = $ws.Range("D35")
$rng.Formula = "=""5.71"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  +0.95%  "
$rng = $ws.Range("D36")
$rng.Formula = "=""5.01"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  +6.34%  "
$ws.Range("E37").Value = "  +0.20%  "
$rng = $ws.Range("D38")
$rng.Formula = "=""4.18"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  -8.36%  "
$rng = $ws.Range("D39")
$rng.Formula = "=""0.0375"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("E41").Value = "  -4.12%  "
$rng = $ws.Range("D42")
$rng.Formula = "=""72.71"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -0.06%  "
$rng = $ws.Range("D43")
$rng.Formula = "=""0.230"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("E44").Value = "  +0.01%  "
$rng = $ws.Range("D45")
$rng.Formula = "=""12.54"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  -7.68%  "
$ws.Range("E46").Value = "  +1.65%  "
$rng = $ws.Range("D47")
$rng.Formula = "=""5.68"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  +1.44%  "
$rng = $ws.Range("D48")
$rng.Formula = "=""108.25"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  -0.50%  "
$rng = $ws.Range("D49")
$rng.Formula = "=""8.87"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("E50").Value = "  -1.15%  "
$rng = $ws.Range("D51")
$rng.Formula = "=""71.98"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  +2.63%  "

$excel.CutCopyMode = 0
